$d = $word.ActiveDocument

# 1) Increase the title font size from 12pt (sz=24) to 13pt (sz=26).
$titleRange = $d.Paragraphs(1).Range
$titleRange.Font.Size = 13
$titleRange.Font.SizeBi = 13

# 2) Merge the two runs of the "Caso algum integrante esteja presente..." paragraph
#    into a single run of text (removing the trailing-space/leading split).
$d.Content.Find.Execute(
    "Caso algum integrante esteja presente na reunião, mas não participe será registrado em ata e computado nos indicadores de participação do projeto.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Caso algum integrante esteja presente na reunião, mas não participe será registrado em ata e computado nos indicadores de participação do projeto.",
    2
) | Out-Null
